# Apply the latest cryptos snapshot (prices in column D, 1h volume % in column E).
# GitHub Actions scheduled update - Thu Apr 18 20:52:41 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.462.39'
$ws.Range("E2").Value = '  +3.99%  '

$ws.Range("D3").Value = '3.062.09'
$ws.Range("E3").Value = '  +2.68%  '

$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '549.07'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.62%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.16'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.61%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.09%  '

$ws.Range("D8").Value = '3.057.24'
$ws.Range("E8").Value = '  +2.78%  '

$ws.Range("E9").Value = '  +1.86%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.54'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.43%  '

$ws.Range("E11").Value = '  +2.39%  '

$ws.Range("E12").Value = '  +2.47%  '

$ws.Range("E13").Value = '  +3.39%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.82'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.56%  '

$ws.Range("D15").Value = '3.563.21'
$ws.Range("E15").Value = '  +2.62%  '

$ws.Range("D16").Value = '63.421.22'
$ws.Range("E16").Value = '  +3.63%  '

$ws.Range("D17").Value = '3.066.18'
$ws.Range("E17").Value = '  +2.28%  '

$ws.Range("E18").Value = '  -1.06%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.76'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.88%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '483.87'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.35%  '

$ws.Range("E21").Value = '  +4.04%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.673'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.39%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.27'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.31%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '80.84'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.90%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.62'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.47%  '

$ws.Range("E26").Value = '  +0.09%  '

$ws.Range("E27").Value = '  +3.06%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.90'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.16%  '

$ws.Range("E29").Value = '  +6.86%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.30%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '26.13'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.94%  '

$ws.Range("E32").Value = '  +0.49%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.45'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +8.40%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.73'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.65%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '55.64'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.29%  '

$ws.Range("E36").Value = '  +2.54%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '465.84'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.61%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0822'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.00%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0397'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.00%  '

$ws.Range("D40").Value = '3.068.56'
$ws.Range("E40").Value = '  -3.92%  '

$ws.Range("E41").Value = '  +1.11%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.26'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.22%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.57'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.40%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '28.06'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.30%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.256'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.29%  '

$ws.Range("E46").Value = '  -0.11%  '

$ws.Range("E47").Value = '  +4.05%  '

$ws.Range("E48").Value = '  +2.23%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '116.34'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.91%  '

$ws.Range("D50").Value = '0.0₃0509'
$ws.Range("E50").Value = '  +3.32%  '

$ws.Range("E51").Value = '  +4.64%  '
